{"js": "// Rename \"GDS Toolkit\" / \"GDS toolkit\" references to \"GOV.UK Toolkit\" / \"GOV.UK toolkit\"\n// throughout the readme, leaving other unrelated \"GDS\" mentions (GDS guidelines, GDS Theme,\n// GDS styling, OutreachGDS* file names, \"GDS Outreach forms\", \"GDS-style help\") untouched.\n// Also adds the extra blank \"Heading 2\" paragraph that appears above the \"Summary\" heading.\n\nconst body = context.document.body;\n\nasync function replacePhrase(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// Title: \"Outreach GDS Toolkit Readme\" -> \"Outreach GOV.UK Toolkit Readme\"\nawait replacePhrase(\"Outreach GDS Toolkit Readme\", \"Outreach GOV.UK Toolkit Readme\");\n\n// Intro paragraph: \"The Outreach gov.uk toolkit allows\" -> \"The Outreach GOV.UK toolkit allows\"\nawait replacePhrase(\"The Outreach gov.uk toolkit allows\", \"The Outreach GOV.UK toolkit allows\");\n\n// \"Installing the GDS toolkit\" heading -> \"Installing the GOV.UK toolkit\"\nawait replacePhrase(\"Installing the GDS toolkit\", \"Installing the GOV.UK toolkit\");\n\n// \"Before installing the GDS Toolkit,\" -> \"Before installing the GOV.UK Toolkit,\"\nawait replacePhrase(\"Before installing the GDS Toolkit,\", \"Before installing the GOV.UK Toolkit,\");\n\n// Four identical installation bullet points: \"... directory of the GDS Toolkit release ...\"\nawait replacePhrase(\" directory of the GDS Toolkit release\", \" directory of the GOV.UK Toolkit release\");\n\n// \"GDS toolkit features\" heading -> \"GOV.UK Toolkit features\"\nawait replacePhrase(\"GDS toolkit features\", \"GOV.UK Toolkit features\");\n\n// \"... configured to use the GDS toolkit will have a ...\" -> \"... GOV.UK Toolkit will have a ...\"\nawait replacePhrase(\"configured to use the GDS toolkit will have a\", \"configured to use the GOV.UK Toolkit will have a\");\n\n// Structural change: an extra blank \"Heading 2\" paragraph was inserted directly above the\n// \"Summary\" heading (it now sits between the existing blank Heading 2 paragraph and \"Summary\").\nconst summaryResults = body.search(\"Summary\", { matchCase: true });\nsummaryResults.load(\"items\");\nawait context.sync();\n\nlet summaryPara = null;\nfor (const item of summaryResults.items) {\n  const p = item.paragraphs.getFirst();\n  p.load(\"style,text\");\n  await context.sync();\n  if (p.style === \"Heading 2\" && p.text === \"Summary\") {\n    summaryPara = p;\n    break;\n  }\n}\n\nif (summaryPara) {\n  const newPara = summaryPara.insertParagraph(\"\", Word.InsertLocation.before);\n  newPara.styleBuiltIn = Word.BuiltInStyleName.heading2;\n  await context.sync();\n}\n", "ps1": "# Rename \"GDS Toolkit\" / \"GDS toolkit\" references to \"GOV.UK Toolkit\" / \"GOV.UK toolkit\"\n# throughout the readme, leaving other unrelated \"GDS\" mentions (GDS guidelines, GDS Theme,\n# GDS styling, OutreachGDS* file names, \"GDS Outreach forms\", \"GDS-style help\") untouched.\n# Also adds the extra blank \"Heading 2\" paragraph that appears above the \"Summary\" heading.\n\nfunction Replace-AllText($doc, $old, $new) {\n    $find = $doc.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n\n$d = $word.ActiveDocument\n\n# Title: \"Outreach GDS Toolkit Readme\" -> \"Outreach GOV.UK Toolkit Readme\"\nReplace-AllText $d \"Outreach GDS Toolkit Readme\" \"Outreach GOV.UK Toolkit Readme\"\n\n# Intro paragraph: \"The Outreach gov.uk toolkit allows\" -> \"The Outreach GOV.UK toolkit allows\"\nReplace-AllText $d \"The Outreach gov.uk toolkit allows\" \"The Outreach GOV.UK toolkit allows\"\n\n# \"Installing the GDS toolkit\" heading -> \"Installing the GOV.UK toolkit\"\nReplace-AllText $d \"Installing the GDS toolkit\" \"Installing the GOV.UK toolkit\"\n\n# \"Before installing the GDS Toolkit,\" -> \"Before installing the GOV.UK Toolkit,\"\nReplace-AllText $d \"Before installing the GDS Toolkit,\" \"Before installing the GOV.UK Toolkit,\"\n\n# Four identical installation bullet points: \"... directory of the GDS Toolkit release ...\"\nReplace-AllText $d \" directory of the GDS Toolkit release\" \" directory of the GOV.UK Toolkit release\"\n\n# \"GDS toolkit features\" heading -> \"GOV.UK Toolkit features\"\nReplace-AllText $d \"GDS toolkit features\" \"GOV.UK Toolkit features\"\n\n# \"... configured to use the GDS toolkit will have a ...\" -> \"... GOV.UK Toolkit will have a ...\"\nReplace-AllText $d \"configured to use the GDS toolkit will have a\" \"configured to use the GOV.UK Toolkit will have a\"\n\n# Structural change: an extra blank \"Heading 2\" paragraph was inserted directly above the\n# \"Summary\" heading (it now sits between the existing blank Heading 2 paragraph and \"Summary\").\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $txt = $p.Range.Text.Trim()\n    if ($txt -eq \"Summary\" -and $p.Style.NameLocal -eq \"Heading 2\") {\n        $target = $p\n        break\n    }\n}\nif ($target -ne $null) {\n    $target.Range.InsertParagraphBefore() | Out-Null\n}\n"}
